$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# --- Fix X24: weapon -> armour (existing row) ---
$ws.Range("X24").Value = "armour"

# --- Append new rows 213-242 ---
# Row 213
$ws.Range("B213").Value = 'Eye For Gold'
$ws.Range("C213").Value = 1
$ws.Range("D213").Value = 'Prince Kalises Golden Plate'
$ws.Range("E213").Value = 'body'
$ws.Range("F213").Value = 'This was worn by the might prince Kalises, he''s dead now. Died from a goblin arrow to the throat. I''ll sell it to you. I stole it from his funeral. Napped it right off his body I did.'
$ws.Range("G213").Value = 'body'
$ws.Range("J213").Value = 190
$ws.Range("K213").Value = 125000
$ws.Range("O213").Value = 0.36
$ws.Range("P213").Value = 0.36
$ws.Range("Q213").Value = 0.36
$ws.Range("R213").Value = 0.36
$ws.Range("S213").Value = 0.36
$ws.Range("U213").Value = 1
$ws.Range("V213").Value = 50
$ws.Range("W213").Value = 100
$ws.Range("X213").Value = 'armour'

# Row 214
$ws.Range("A214").Value = 'Balanced Energies'
$ws.Range("C214").Value = 1
$ws.Range("D214").Value = 'Long Silver Sword'
$ws.Range("E214").Value = 'weapon'
$ws.Range("F214").Value = 'A simple long sword made of the finest silver.'
$ws.Range("H214").Value = 8
$ws.Range("K214").Value = 50
$ws.Range("O214").Value = 0.05
$ws.Range("P214").Value = 0.05
$ws.Range("Q214").Value = 0.05
$ws.Range("R214").Value = 0.05
$ws.Range("S214").Value = 0.05
$ws.Range("U214").Value = 1
$ws.Range("V214").Value = 3
$ws.Range("W214").Value = 8
$ws.Range("X214").Value = 'weapon'

# Row 215
$ws.Range("B215").Value = 'Sinister Dance'
$ws.Range("C215").Value = 1
$ws.Range("D215").Value = 'Thorn Rose Plate'
$ws.Range("E215").Value = 'body'
$ws.Range("F215").Value = 'Covered in metal thorns, this plate is painted white with the images of roses and vibrant thorns.'
$ws.Range("G215").Value = 'body'
$ws.Range("J215").Value = 170
$ws.Range("K215").Value = 48000
$ws.Range("O215").Value = 0.33
$ws.Range("P215").Value = 0.33
$ws.Range("Q215").Value = 0.33
$ws.Range("R215").Value = 0.33
$ws.Range("S215").Value = 0.33
$ws.Range("U215").Value = 1
$ws.Range("V215").Value = 43
$ws.Range("W215").Value = 78
$ws.Range("X215").Value = 'armour'

# Row 216
$ws.Range("A216").Value = 'Blood Lust'
$ws.Range("C216").Value = 1
$ws.Range("D216").Value = 'Long Silver Sword'
$ws.Range("E216").Value = 'weapon'
$ws.Range("F216").Value = 'A simple long sword made of the finest silver.'
$ws.Range("H216").Value = 8
$ws.Range("K216").Value = 50
$ws.Range("O216").Value = 0.05
$ws.Range("P216").Value = 0.05
$ws.Range("Q216").Value = 0.05
$ws.Range("R216").Value = 0.05
$ws.Range("S216").Value = 0.05
$ws.Range("U216").Value = 1
$ws.Range("V216").Value = 3
$ws.Range("W216").Value = 8
$ws.Range("X216").Value = 'weapon'

# Row 217
$ws.Range("A217").Value = 'Blood Lust'
$ws.Range("C217").Value = 1
$ws.Range("D217").Value = 'Broken Dagger'
$ws.Range("E217").Value = 'weapon'
$ws.Range("F217").Value = 'A simple, shattered broken dagger. It''s at least still sharp.'
$ws.Range("H217").Value = 4
$ws.Range("K217").Value = 10
$ws.Range("U217").Value = 1
$ws.Range("V217").Value = 1
$ws.Range("W217").Value = 5
$ws.Range("X217").Value = 'weapon'

# Row 218
$ws.Range("A218").Value = 'Fighters Strength'
$ws.Range("C218").Value = 1
$ws.Range("D218").Value = 'Broken Dagger'
$ws.Range("E218").Value = 'weapon'
$ws.Range("F218").Value = 'A simple, shattered broken dagger. It''s at least still sharp.'
$ws.Range("H218").Value = 4
$ws.Range("K218").Value = 10
$ws.Range("U218").Value = 1
$ws.Range("V218").Value = 1
$ws.Range("W218").Value = 5
$ws.Range("X218").Value = 'weapon'

# Row 219
$ws.Range("A219").Value = 'Archers Bane'
$ws.Range("C219").Value = 1
$ws.Range("D219").Value = 'Broken Dagger'
$ws.Range("E219").Value = 'weapon'
$ws.Range("F219").Value = 'A simple, shattered broken dagger. It''s at least still sharp.'
$ws.Range("H219").Value = 4
$ws.Range("K219").Value = 10
$ws.Range("U219").Value = 1
$ws.Range("V219").Value = 1
$ws.Range("W219").Value = 5
$ws.Range("X219").Value = 'weapon'

# Row 220
$ws.Range("B220").Value = 'Dancers Moves'
$ws.Range("C220").Value = 1
$ws.Range("D220").Value = 'Dark Steel Breast Plate'
$ws.Range("E220").Value = 'body'
$ws.Range("F220").Value = 'Made from the rarest of steel, dark steel. No one know whats stains steel this dark when it''s made. It is stronger and more durable the other forms of steel.'
$ws.Range("G220").Value = 'body'
$ws.Range("J220").Value = 40
$ws.Range("K220").Value = 750
$ws.Range("O220").Value = 0.18
$ws.Range("P220").Value = 0.18
$ws.Range("Q220").Value = 0.18
$ws.Range("R220").Value = 0.18
$ws.Range("S220").Value = 0.18
$ws.Range("U220").Value = 1
$ws.Range("V220").Value = 18
$ws.Range("W220").Value = 36
$ws.Range("X220").Value = 'armour'

# Row 221
$ws.Range("A221").Value = 'Archbishops Prayer'
$ws.Range("B221").Value = 'Desert Winds'
$ws.Range("C221").Value = 1
$ws.Range("D221").Value = 'Tiger Tooth Knife'
$ws.Range("E221").Value = 'weapon'
$ws.Range("F221").Value = 'Made from a tigers tooth. This knife is easy to conceal and easy to use.'
$ws.Range("H221").Value = 450
$ws.Range("K221").Value = 1525000
$ws.Range("O221").Value = 0.37
$ws.Range("P221").Value = 0.37
$ws.Range("Q221").Value = 0.37
$ws.Range("R221").Value = 0.37
$ws.Range("S221").Value = 0.37
$ws.Range("U221").Value = 1
$ws.Range("V221").Value = 70
$ws.Range("W221").Value = 150
$ws.Range("X221").Value = 'weapon'

# Row 222
$ws.Range("B222").Value = 'Eye For Gold'
$ws.Range("C222").Value = 1
$ws.Range("D222").Value = 'Sage and Thread'
$ws.Range("E222").Value = 'spell-healing'
$ws.Range("F222").Value = 'Through the use of sage and magical thread, you''ll cure your self in no time.'
$ws.Range("I222").Value = 260
$ws.Range("K222").Value = 320000
$ws.Range("U222").Value = 1
$ws.Range("V222").Value = 55
$ws.Range("W222").Value = 110
$ws.Range("X222").Value = 'spell'

# Row 223
$ws.Range("B223").Value = 'Enchantress Luck'
$ws.Range("C223").Value = 1
$ws.Range("D223").Value = 'Broken Wooden Shield'
$ws.Range("E223").Value = 'shield'
$ws.Range("F223").Value = 'This use to be a wooden shield, now its nothing more then a plank of wood with some metal. Have fun.'
$ws.Range("J223").Value = 1
$ws.Range("K223").Value = 5
$ws.Range("U223").Value = 1
$ws.Range("V223").Value = 1
$ws.Range("W223").Value = 5
$ws.Range("X223").Value = 'armour'

# Row 224
$ws.Range("B224").Value = 'Enchantress Luck'
$ws.Range("C224").Value = 1
$ws.Range("D224").Value = 'Oak Shield'
$ws.Range("E224").Value = 'shield'
$ws.Range("F224").Value = 'Made completely of oak and strong. This was some good carpentry.'
$ws.Range("J224").Value = 3
$ws.Range("K224").Value = 60
$ws.Range("O224").Value = 0.02
$ws.Range("P224").Value = 0.02
$ws.Range("Q224").Value = 0.02
$ws.Range("R224").Value = 0.02
$ws.Range("S224").Value = 0.02
$ws.Range("U224").Value = 1
$ws.Range("V224").Value = 3
$ws.Range("W224").Value = 8
$ws.Range("X224").Value = 'armour'

# Row 225
$ws.Range("A225").Value = 'Mages Inspiration'
$ws.Range("C225").Value = 1
$ws.Range("D225").Value = 'Necromancers Gloves'
$ws.Range("E225").Value = 'gloves'
$ws.Range("F225").Value = 'Worn by necromancers who raise the dead, these gloves have the blood of their enemies all over them.'
$ws.Range("J225").Value = 150
$ws.Range("K225").Value = 360000
$ws.Range("O225").Value = 0.18
$ws.Range("P225").Value = 0.18
$ws.Range("Q225").Value = 0.18
$ws.Range("R225").Value = 0.18
$ws.Range("S225").Value = 0.18
$ws.Range("U225").Value = 1
$ws.Range("V225").Value = 55
$ws.Range("W225").Value = 110
$ws.Range("X225").Value = 'armour'

# Row 226
$ws.Range("A226").Value = 'Soldiers Resilance'
$ws.Range("C226").Value = 1
$ws.Range("D226").Value = 'Tin Helmet'
$ws.Range("E226").Value = 'helmet'
$ws.Range("F226").Value = 'Simple, not very sturdy, but simple.'
$ws.Range("G226").Value = 'helmet'
$ws.Range("J226").Value = 3
$ws.Range("K226").Value = 55
$ws.Range("U226").Value = 1
$ws.Range("V226").Value = 3
$ws.Range("W226").Value = 8
$ws.Range("X226").Value = 'armour'

# Row 227
$ws.Range("B227").Value = 'Desert Winds'
$ws.Range("C227").Value = 1
$ws.Range("D227").Value = 'Life Stealing'
$ws.Range("E227").Value = 'spell-damage'
$ws.Range("F227").Value = 'Steal the life from the enemy'
$ws.Range("H227").Value = 200
$ws.Range("K227").Value = 1250000
$ws.Range("U227").Value = 1
$ws.Range("V227").Value = 70
$ws.Range("W227").Value = 150
$ws.Range("X227").Value = 'spell'

# Row 228
$ws.Range("A228").Value = 'Fighters Strength'
$ws.Range("C228").Value = 1
$ws.Range("D228").Value = 'Long Silver Sword'
$ws.Range("E228").Value = 'weapon'
$ws.Range("F228").Value = 'A simple long sword made of the finest silver.'
$ws.Range("H228").Value = 8
$ws.Range("K228").Value = 50
$ws.Range("O228").Value = 0.05
$ws.Range("P228").Value = 0.05
$ws.Range("Q228").Value = 0.05
$ws.Range("R228").Value = 0.05
$ws.Range("S228").Value = 0.05
$ws.Range("U228").Value = 1
$ws.Range("V228").Value = 3
$ws.Range("W228").Value = 8
$ws.Range("X228").Value = 'weapon'

# Row 229
$ws.Range("B229").Value = 'Deaths Accuracy'
$ws.Range("C229").Value = 1
$ws.Range("D229").Value = 'Flail Of Hell'
$ws.Range("E229").Value = 'weapon'
$ws.Range("F229").Value = 'Created in the pits of hell this flail has the head of the mace as a ball of fire.'
$ws.Range("H229").Value = 380
$ws.Range("K229").Value = 750000
$ws.Range("O229").Value = 0.35
$ws.Range("P229").Value = 0.35
$ws.Range("Q229").Value = 0.35
$ws.Range("R229").Value = 0.35
$ws.Range("S229").Value = 0.35
$ws.Range("U229").Value = 1
$ws.Range("V229").Value = 60
$ws.Range("W229").Value = 130
$ws.Range("X229").Value = 'weapon'

# Row 230
$ws.Range("A230").Value = 'Soldiers Resilance'
$ws.Range("C230").Value = 1
$ws.Range("D230").Value = 'Flower Of Roots Club'
$ws.Range("E230").Value = 'weapon'
$ws.Range("F230").Value = 'Made of roots that flower, this club sure is unique.'
$ws.Range("H230").Value = 280
$ws.Range("K230").Value = 125000
$ws.Range("O230").Value = 0.29
$ws.Range("P230").Value = 0.29
$ws.Range("Q230").Value = 0.29
$ws.Range("R230").Value = 0.29
$ws.Range("S230").Value = 0.29
$ws.Range("U230").Value = 1
$ws.Range("V230").Value = 50
$ws.Range("W230").Value = 100
$ws.Range("X230").Value = 'weapon'

# Row 231
$ws.Range("B231").Value = 'Armour Smiths Hopes'
$ws.Range("C231").Value = 1
$ws.Range("D231").Value = 'Dark Steel Breast Plate'
$ws.Range("E231").Value = 'body'
$ws.Range("F231").Value = 'Made from the rarest of steel, dark steel. No one know whats stains steel this dark when it''s made. It is stronger and more durable the other forms of steel.'
$ws.Range("G231").Value = 'body'
$ws.Range("J231").Value = 40
$ws.Range("K231").Value = 750
$ws.Range("O231").Value = 0.18
$ws.Range("P231").Value = 0.18
$ws.Range("Q231").Value = 0.18
$ws.Range("R231").Value = 0.18
$ws.Range("S231").Value = 0.18
$ws.Range("U231").Value = 1
$ws.Range("V231").Value = 18
$ws.Range("W231").Value = 36
$ws.Range("X231").Value = 'armour'

# Row 232
$ws.Range("A232").Value = 'Soldiers Resilance'
$ws.Range("C232").Value = 1
$ws.Range("D232").Value = 'Dark Fairy Leaf Plate'
$ws.Range("E232").Value = 'body'
$ws.Range("F232").Value = 'Somewhere in the woods is a home of dark fairies. They make a plate that is string and durable.'
$ws.Range("G232").Value = 'body'
$ws.Range("J232").Value = 60
$ws.Range("K232").Value = 1760
$ws.Range("O232").Value = 0.2
$ws.Range("P232").Value = 0.2
$ws.Range("Q232").Value = 0.2
$ws.Range("R232").Value = 0.2
$ws.Range("S232").Value = 0.2
$ws.Range("U232").Value = 1
$ws.Range("V232").Value = 24
$ws.Range("W232").Value = 40
$ws.Range("X232").Value = 'armour'

# Row 233
$ws.Range("B233").Value = 'Sinister Dance'
$ws.Range("C233").Value = 1
$ws.Range("D233").Value = 'Enchantress Gloves'
$ws.Range("E233").Value = 'gloves'
$ws.Range("F233").Value = 'Simple gloves made and worn by a now dead enchantress. I stole these off her body.'
$ws.Range("G233").Value = 'gloves'
$ws.Range("J233").Value = 10
$ws.Range("K233").Value = 200
$ws.Range("O233").Value = 0.03
$ws.Range("P233").Value = 0.03
$ws.Range("Q233").Value = 0.03
$ws.Range("R233").Value = 0.03
$ws.Range("S233").Value = 0.03
$ws.Range("U233").Value = 1
$ws.Range("V233").Value = 10
$ws.Range("W233").Value = 20
$ws.Range("X233").Value = 'armour'

# Row 234
$ws.Range("B234").Value = 'Eye For Gold'
$ws.Range("C234").Value = 1
$ws.Range("D234").Value = 'Crystal Ring'
$ws.Range("E234").Value = 'ring'
$ws.Range("F234").Value = 'Made of crystals this ring is vibrating with energy'
$ws.Range("H234").Value = 55
$ws.Range("K234").Value = 13400
$ws.Range("L234").Value = 0.16
$ws.Range("M234").Value = 0.16
$ws.Range("N234").Value = 0.16
$ws.Range("U234").Value = 1
$ws.Range("V234").Value = 36
$ws.Range("W234").Value = 70
$ws.Range("X234").Value = 'ring'

# Row 235
$ws.Range("A235").Value = 'Blood Lust'
$ws.Range("C235").Value = 1
$ws.Range("D235").Value = 'Dark Steel Breast Plate'
$ws.Range("E235").Value = 'body'
$ws.Range("F235").Value = 'Made from the rarest of steel, dark steel. No one know whats stains steel this dark when it''s made. It is stronger and more durable the other forms of steel.'
$ws.Range("G235").Value = 'body'
$ws.Range("J235").Value = 40
$ws.Range("K235").Value = 750
$ws.Range("O235").Value = 0.18
$ws.Range("P235").Value = 0.18
$ws.Range("Q235").Value = 0.18
$ws.Range("R235").Value = 0.18
$ws.Range("S235").Value = 0.18
$ws.Range("U235").Value = 1
$ws.Range("V235").Value = 18
$ws.Range("W235").Value = 36
$ws.Range("X235").Value = 'armour'

# Row 236
$ws.Range("B236").Value = 'Blacksmiths Heart'
$ws.Range("C236").Value = 1
$ws.Range("D236").Value = 'Dark Steel Breast Plate'
$ws.Range("E236").Value = 'body'
$ws.Range("F236").Value = 'Made from the rarest of steel, dark steel. No one know whats stains steel this dark when it''s made. It is stronger and more durable the other forms of steel.'
$ws.Range("G236").Value = 'body'
$ws.Range("J236").Value = 40
$ws.Range("K236").Value = 750
$ws.Range("O236").Value = 0.18
$ws.Range("P236").Value = 0.18
$ws.Range("Q236").Value = 0.18
$ws.Range("R236").Value = 0.18
$ws.Range("S236").Value = 0.18
$ws.Range("U236").Value = 1
$ws.Range("V236").Value = 18
$ws.Range("W236").Value = 36
$ws.Range("X236").Value = 'armour'

# Row 237
$ws.Range("A237").Value = 'Balanced Energies'
$ws.Range("C237").Value = 1
$ws.Range("D237").Value = 'Iron Clad Robes'
$ws.Range("E237").Value = 'body'
$ws.Range("F237").Value = 'Not sure how this works, there is iron all over this robe, yet some how it works.'
$ws.Range("G237").Value = 'body'
$ws.Range("J237").Value = 230
$ws.Range("K237").Value = 450000
$ws.Range("O237").Value = 0.4
$ws.Range("P237").Value = 0.4
$ws.Range("Q237").Value = 0.4
$ws.Range("R237").Value = 0.4
$ws.Range("S237").Value = 0.4
$ws.Range("U237").Value = 1
$ws.Range("V237").Value = 55
$ws.Range("W237").Value = 110
$ws.Range("X237").Value = 'armour'

# Row 238
$ws.Range("A238").Value = 'Hawk Eye'
$ws.Range("C238").Value = 1
$ws.Range("D238").Value = 'Dark Fairy Leaf Plate'
$ws.Range("E238").Value = 'body'
$ws.Range("F238").Value = 'Somewhere in the woods is a home of dark fairies. They make a plate that is string and durable.'
$ws.Range("G238").Value = 'body'
$ws.Range("J238").Value = 60
$ws.Range("K238").Value = 1760
$ws.Range("O238").Value = 0.2
$ws.Range("P238").Value = 0.2
$ws.Range("Q238").Value = 0.2
$ws.Range("R238").Value = 0.2
$ws.Range("S238").Value = 0.2
$ws.Range("U238").Value = 1
$ws.Range("V238").Value = 24
$ws.Range("W238").Value = 40
$ws.Range("X238").Value = 'armour'

# Row 239
$ws.Range("A239").Value = 'Clerics Blessing'
$ws.Range("C239").Value = 1
$ws.Range("D239").Value = 'Angels Vengeance'
$ws.Range("E239").Value = 'spell-damage'
$ws.Range("F239").Value = 'Angels have a vengeance that they can take out on the enemies of the faithful.'
$ws.Range("H239").Value = 120
$ws.Range("K239").Value = 55000
$ws.Range("U239").Value = 1
$ws.Range("V239").Value = 43
$ws.Range("W239").Value = 78
$ws.Range("X239").Value = 'spell'

# Row 240
$ws.Range("A240").Value = 'Balanced Energies'
$ws.Range("C240").Value = 1
$ws.Range("D240").Value = 'Goblin Treasure Ring'
$ws.Range("E240").Value = 'ring'
$ws.Range("F240").Value = 'Goblins got stories about treasure rings bringing great treasure.'
$ws.Range("H240").Value = 90
$ws.Range("K240").Value = 750000
$ws.Range("L240").Value = 0.29
$ws.Range("M240").Value = 0.29
$ws.Range("N240").Value = 0.29
$ws.Range("U240").Value = 1
$ws.Range("V240").Value = 60
$ws.Range("W240").Value = 130
$ws.Range("X240").Value = 'ring'

# Row 241
$ws.Range("A241").Value = 'Fighters Strength'
$ws.Range("C241").Value = 1
$ws.Range("D241").Value = 'Empaths Gloves'
$ws.Range("E241").Value = 'gloves'
$ws.Range("F241").Value = 'Worn by those who are sensitive to others emotions and feelings, these gloves may give you some empathy.'
$ws.Range("G241").Value = 'gloves'
$ws.Range("J241").Value = 200
$ws.Range("K241").Value = 800000
$ws.Range("O241").Value = 0.23
$ws.Range("P241").Value = 0.23
$ws.Range("Q241").Value = 0.23
$ws.Range("R241").Value = 0.23
$ws.Range("S241").Value = 0.23
$ws.Range("U241").Value = 1
$ws.Range("V241").Value = 70
$ws.Range("W241").Value = 180
$ws.Range("X241").Value = 'armour'

# Row 242
$ws.Range("A242").Value = 'Fighters Strength'
$ws.Range("C242").Value = 1
$ws.Range("D242").Value = 'Holy Ring'
$ws.Range("E242").Value = 'ring'
$ws.Range("F242").Value = 'Blessed by the high priests across the sea.'
$ws.Range("H242").Value = 80
$ws.Range("K242").Value = 346000
$ws.Range("L242").Value = 0.27
$ws.Range("M242").Value = 0.27
$ws.Range("N242").Value = 0.27
$ws.Range("U242").Value = 1
$ws.Range("V242").Value = 55
$ws.Range("W242").Value = 110
$ws.Range("X242").Value = 'ring'
